# "release chapter highly cited"
#
# Nudges the position/size of a handful of shapes inside the single
# group shape on slide 1 (the world-map plot):
#   - rc3  (background rect)      : width shrinks by 1 EMU
#   - tx8  ("Asia" label)         : re-positioned slightly
#   - tx9  ("Europe" label)       : re-positioned slightly
#   - tx10 ("North America" label): re-positioned slightly
#   - tx11 ("South America" label): re-positioned slightly
#
# Shape.Left/Top/Width/Height are exposed in points (1 pt = 12700 EMU)
# and are backed by 32-bit floats, so literals below are chosen to
# round-trip to the exact target EMU values instead of the naive
# EMU/12700 division.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$g = $s.Shapes.Item(1)

# rc3 : ext cx 6400800 -> 6400799 (cy unchanged)
$rc3 = $g.GroupItems.Item(1)
$rc3.Width = 503.99993

# tx8 ("Asia") : off 4779076,2841724 -> 4779444,2842688
$tx8 = $g.GroupItems.Item(6)
$tx8.Left = 376.3342
$tx8.Top = 223.83371

# tx9 ("Europe") : off 4381873,4470147 -> 4382594,4472716
$tx9 = $g.GroupItems.Item(7)
$tx9.Left = 345.08614
$tx9.Top = 352.1824

# tx10 ("North America") : off 5832638,3346765 -> 5834563,3350631
$tx10 = $g.GroupItems.Item(8)
$tx10.Left = 459.41442
$tx10.Top = 263.82921

# tx11 ("South America") : off 5621140,2844509 -> 5621431,2850374
$tx11 = $g.GroupItems.Item(9)
$tx11.Left = 442.6324
$tx11.Top = 224.4389
